$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1454.6666
$ws.Range("I32").Value = 1559.3334
$ws.Range("J32").Value = 1350
$ws.Range("K32").Value = 1559.3334
$ws.Range("L32").Value = 1350
$ws.Range("M32").Value = -1233.3334
$ws.Range("N32").Value = -2002
$ws.Range("H129").Value = 2189.375
$ws.Range("J129").Value = 2496.1667
$ws.Range("L129").Value = 7488.500100000001
$ws.Range("N129").Value = -17488.5001
$ws.Range("H141").Value = 2400
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 555.8570999999999
$ws.Range("I2").Value = 578.4
$ws.Range("J2").Value = 499.5
$ws.Range("K2").Value = 578.4
$ws.Range("L2").Value = 499.5
$ws.Range("M2").Value = -465.4
$ws.Range("N2").Value = -725.5
$ws.Range("H74").Value = 3999
$ws.Range("I74").Value = 3999
$ws.Range("K74").Value = 3999
$ws.Range("M74").Value = -3125
$ws.Range("H77").Value = 3999
$ws.Range("I77").Value = 3999
$ws.Range("K77").Value = 19995
$ws.Range("M77").Value = -15627
$ws.Range("H116").Value = 555.8570999999999
$ws.Range("I116").Value = 578.4
$ws.Range("J116").Value = 499.5
$ws.Range("K116").Value = 578.4
$ws.Range("L116").Value = 499.5
$ws.Range("M116").Value = 1715.6
$ws.Range("N116").Value = -5087.5
$ws.Range("H132").Value = 1714.75
$ws.Range("I132").Value = 1069.4138
$ws.Range("J132").Value = 3416.0908
$ws.Range("K132").Value = 3208.2414
$ws.Range("L132").Value = 10248.2724
$ws.Range("M132").Value = -678.2413999999999
$ws.Range("N132").Value = -15308.2724

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 555.8570999999999
$ws.Range("I3").Value = 578.4
$ws.Range("J3").Value = 499.5
$ws.Range("K3").Value = 578.4
$ws.Range("L3").Value = 499.5
$ws.Range("M3").Value = -464.4
$ws.Range("N3").Value = -727.5
$ws.Range("H86").Value = 13832.333
$ws.Range("I86").Value = 2333.3333
$ws.Range("K86").Value = 2333.3333
$ws.Range("M86").Value = -1210.3333
$ws.Range("H89").Value = 13832.333
$ws.Range("I89").Value = 2333.3333
$ws.Range("K89").Value = 11666.6665
$ws.Range("M89").Value = -6050.666499999999
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H107").Value = 700
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 700
$ws.Range("M107").Value = 1220
$ws.Range("H134").Value = 6631.1904
$ws.Range("I134").Value = 1018.2143
$ws.Range("J134").Value = 17857.143
$ws.Range("K134").Value = 3054.6429
$ws.Range("L134").Value = 53571.429
$ws.Range("M134").Value = -519.6428999999998
$ws.Range("N134").Value = -58641.429

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2833.3333
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2833.3333
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 746.125
$ws.Range("I107").Value = 661.5
$ws.Range("K107").Value = 661.5
$ws.Range("M107").Value = 1258.5
$ws.Range("H132").Value = 2621.4443
$ws.Range("I132").Value = 1942.1428
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5826.428400000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3296.428400000001
$ws.Range("N132").Value = -20057

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3331
$ws.Range("I3").Value = 2222
$ws.Range("J3").Value = 4440
$ws.Range("K3").Value = 6666
$ws.Range("L3").Value = 13320
$ws.Range("M3").Value = -6554
$ws.Range("N3").Value = -13544
$ws.Range("H18").Value = 1202.5
$ws.Range("I18").Value = 1202.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 3607.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3438.5
$ws.Range("N18").ClearContents()
$ws.Range("H109").Value = 4572.6665
$ws.Range("I109").Value = 2756.5
$ws.Range("K109").Value = 8269.5
$ws.Range("M109").Value = -7229.5
$ws.Range("H111").Value = 2075.6667
$ws.Range("I111").Value = 863.5
$ws.Range("J111").Value = 4500
$ws.Range("K111").Value = 2590.5
$ws.Range("L111").Value = 13500
$ws.Range("M111").Value = 476.5
$ws.Range("N111").Value = -19634
$ws.Range("H114").Value = 2189.8
$ws.Range("I114").Value = 2129.6
$ws.Range("J114").Value = 2250
$ws.Range("K114").Value = 6388.799999999999
$ws.Range("L114").Value = 6750
$ws.Range("M114").Value = -3134.799999999999
$ws.Range("N114").Value = -13258
$ws.Range("H121").Value = 846.6667
$ws.Range("I121").Value = 495
$ws.Range("K121").Value = 1485
$ws.Range("M121").Value = -175
$ws.Range("H123").Value = 5000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H134").Value = 4076.75
$ws.Range("I134").Value = 4076.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 12230.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7160.25
$ws.Range("N134").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 382.5
$ws.Range("I97").Value = 384.5
$ws.Range("J97").Value = 376.5
$ws.Range("K97").Value = 384.5
$ws.Range("L97").Value = 376.5
$ws.Range("M97").Value = 111.5
$ws.Range("N97").Value = -1368.5
$ws.Range("H102").Value = 2501.375
$ws.Range("I102").Value = 1668.6666
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 1668.6666
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = -46.66660000000002
$ws.Range("N102").Value = -8243.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 966.6667
$ws.Range("I16").Value = 966.6667
$ws.Range("K16").Value = 966.6667
$ws.Range("M16").Value = -796.6667
$ws.Range("H22").Value = 1686.125
$ws.Range("J22").Value = 1740
$ws.Range("L22").Value = 1740
$ws.Range("N22").Value = -2330
$ws.Range("H27").Value = 1686.125
$ws.Range("J27").Value = 1740
$ws.Range("L27").Value = 1740
$ws.Range("N27").Value = -1954

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1241.7858
$ws.Range("I113").Value = 1497.7142
$ws.Range("J113").Value = 985.8570999999999
$ws.Range("K113").Value = 4493.142599999999
$ws.Range("L113").Value = 2957.5713
$ws.Range("M113").Value = -2323.142599999999
$ws.Range("N113").Value = -7297.5713
$ws.Range("H132").Value = 5287.9287
$ws.Range("I132").Value = 2521.5
$ws.Range("K132").Value = 7564.5
$ws.Range("M132").Value = -5034.5
$ws.Range("H136").Value = 1763.0667
$ws.Range("I136").Value = 1746.1428
$ws.Range("K136").Value = 5238.428400000001
$ws.Range("M136").Value = -2688.428400000001
